# The deck currently has its Design/Theme set to the custom "Integral"
# theme (green palette). The target revision swaps the presentation
# back to the default "Office Theme" palette (blue palette) - i.e. the
# 12 theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# stored in the active theme part are reset to their stock Office
# Theme values.
#
# Helper: convert a "RRGGBB" hex string into the BGR-packed integer
# that PowerPoint's RGB-valued COM properties expect.
function Convert-HexToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Stock "Office Theme" 12-slot color scheme, in the same order as
# ThemeColorScheme.Item(1..12): dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink.
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $tcs.Item($i).RGB = Convert-HexToBgr $officeThemeHex[$i - 1]
}
